# Updated symbol list on Wed Dec 14 17:57:06 UTC 2022 with GitHub Actions
# Refreshes the "Price" (column D) values for the coin rows whose quotes changed.
#
# The Price cells are stored as text (not numbers), so for every touched cell we
# first force the cell's number format to Text ("@") and then assign the new
# value as a string. This prevents Excel from auto-coercing the numeric-looking
# text into a floating point number (which would both change the cell's type
# and could introduce binary rounding artifacts), keeping the cell a text value
# exactly as it appears in the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellAddress,
        [string]$NewValue
    )
    $range = $ws.Range($CellAddress)
    $range.NumberFormat = "@"
    $range.Value = $NewValue
}

Set-TextValue "D2"  "272.24"
Set-TextValue "D4"  "6.358"
Set-TextValue "D5"  "0.06296"
Set-TextValue "D8"  "1.403"
Set-TextValue "D9"  "0.8349"
Set-TextValue "D10" "0.1630"
Set-TextValue "D11" "0.08414"
Set-TextValue "D12" "0.03472"
Set-TextValue "D13" "0.03137"
Set-TextValue "D14" "0.09320"
Set-TextValue "D15" "3.948"
Set-TextValue "D16" "0.001709"
Set-TextValue "D17" "0.04863"
Set-TextValue "D18" "0.006279"
Set-TextValue "D19" "0.005478"
Set-TextValue "D20" "0.001090"
Set-TextValue "D22" "3.732"
Set-TextValue "D23" "2.369"
Set-TextValue "D24" "0.01386"
Set-TextValue "D26" "0.1217"
Set-TextValue "D40" "0.04685"
Set-TextValue "D41" "0.006905"
Set-TextValue "D42" "0.1178"
Set-TextValue "D43" "0.003351"
Set-TextValue "D45" "0.00006273"
Set-TextValue "D48" "0.1172"
